$d = $word.ActiveDocument

# 1. Delete the first paragraph "Introduction" (the plain, non-heading one).
#    It is the very first paragraph of the document.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = ""
$p1.Range.Delete()

# 2. Change the text of the "Introduction" Heading1 paragraph to "Assumptions".
$d.Content.Find.Execute("Introduction", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Assumptions", 2)

# 3. Replace the long "Assumption: ..." paragraph text with the new Customs
#    Management Module text.
$old3 = "Assumption: The functional specification document will specify the details for generating a customer invoice form for Canada bound by the profit centers and profit center groups which must be already created and activated prior to execution of the upload."
$new3 = "The functional specification document would assume the setup of a Customs Management Module to generate a custom invoice form tailored to the specific needs of Canadian customers. It would set up a document type to trigger the Adobe Forms in BTP and define the needed information and fields to include in the custom invoice form."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new3, 2)

# 4. Replace the closing paragraph's text with the new VF01/VF02/VF03 text.
$old4 = "Assuming from the given requirements, the functional specification document will include the details of entering profit center information in hierarchical format, data cleansing and signing off by business, and details of error management process for rectifying the errors in uploaded data using an error log file and SM35 transaction code."
$new4 = "The functional specification document will specify the details for generating a custom invoice form in Adobe Forms BTP for Canadian customers. It will outline the information and details, creating master data for customization and creating invoice documents in VF01/VF02/VF03 in SAP for receipt and re-print in VF31. The invoice will be issued either through Print, Email or FAX."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new4, 2)
